$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 759, shifting existing rows 759-807 down to 761-809
$ws.Rows.Item(759).Resize(2).Insert()

# New row 759: Red Globe, Provincia del Elquí
$r = 759
$ws.Cells.Item($r, 1).Value = 5
$ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($r, 3).Value = "Maule"
$ws.Cells.Item($r, 4).Value = 45265
$ws.Cells.Item($r, 5).Value = 7
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100109
$ws.Cells.Item($r, 8).Value = "Uva"
$ws.Cells.Item($r, 9).Value = 100109001
$ws.Cells.Item($r, 10).Value = "Uva"
$ws.Cells.Item($r, 11).Value = "Red Globe"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 16000
$ws.Cells.Item($r, 15).Value = 16000
$ws.Cells.Item($r, 16).Value = 16000
$ws.Cells.Item($r, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item($r, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($r, 19).Value = 2000
$ws.Cells.Item($r, 20).Value = 8

# New row 760: Superior Seedless, Provincia del Elquí
$r = 760
$ws.Cells.Item($r, 1).Value = 5
$ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($r, 3).Value = "Maule"
$ws.Cells.Item($r, 4).Value = 45265
$ws.Cells.Item($r, 5).Value = 7
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100109
$ws.Cells.Item($r, 8).Value = "Uva"
$ws.Cells.Item($r, 9).Value = 100109001
$ws.Cells.Item($r, 10).Value = "Uva"
$ws.Cells.Item($r, 11).Value = "Superior Seedless"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 360
$ws.Cells.Item($r, 14).Value = 15000
$ws.Cells.Item($r, 15).Value = 15000
$ws.Cells.Item($r, 16).Value = 15000
$ws.Cells.Item($r, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item($r, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($r, 19).Value = 1875
$ws.Cells.Item($r, 20).Value = 8
